$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.1059617549180984,0.9682161808013916,0.0132368877530098,0.998595654964447),
    @(0.03349165990948677,0.9897164106369019,0.006149070337414742,0.9988117218017578),
    @(0.02115354500710964,0.993182361125946,0.002184608718380332,0.9990277886390686),
    @(0.01290764007717371,0.9960008263587952,0.002385798608884215,0.9995678663253784),
    @(0.01142254751175642,0.9961341619491577,0.00174883846193552,0.9996759295463562),
    @(0.01090620085597038,0.996400773525238,0.001958071952685714,0.9996759295463562),
    @(0.01061406265944242,0.9964388608932495,0.0009690545848570764,0.9996759295463562),
    @(0.01164072193205357,0.9958294630050659,0.0005478461389429867,0.9997839331626892),
    @(0.01080966275185347,0.9959818124771118,0.001090605393983424,0.9996759295463562),
    @(0.01117768604308367,0.9961150884628296,0.0004921670770272613,0.9997839331626892),
    @(0.01024078670889139,0.996305525302887,0.0004142569960094988,0.999891996383667),
    @(0.01014704816043377,0.996305525302887,0.0006430239300243556,0.9997839331626892),
    @(0.0106438510119915,0.9959437251091003,0.0007191959302872419,0.9995678663253784),
    @(0.0102799516171217,0.9962674379348755,0.000441723212134093,0.9997839331626892),
    @(0.01002911385148764,0.996305525302887,0.0002851988829206675,0.999891996383667),
    @(0.008964367210865021,0.9960770010948181,0.0006866762414574623,0.9996759295463562),
    @(0.008045049384236336,0.9964578747749329,0.0004027200629934669,0.9997839331626892),
    @(0.008178435266017914,0.9960579872131348,0.0005334122688509524,0.9996759295463562),
    @(0.0086591811850667,0.9959246516227722,0.0007029973203316331,0.9996759295463562),
    @(0.008197599090635777,0.9961150884628296,0.0009974768618121743,0.9996759295463562),
    @(0.009743847884237766,0.9955628514289856,0.0009251992451027036,0.9996759295463562),
    @(0.008353312499821186,0.9959627389907837,0.0006093078409321606,0.9996759295463562),
    @(0.006967070046812296,0.9967625737190247,0.000425924634328112,0.9996759295463562),
    @(0.007844255305826664,0.9961150884628296,0.002170724561437964,0.9994598627090454),
    @(0.009217933751642704,0.9959818124771118,0.0008829801226966083,0.9996759295463562),
    @(0.006879780441522598,0.9969720840454102,0.0004762225726153702,0.9996759295463562),
    @(0.008077841252088547,0.9962865114212036,0.0002237200533272699,1),
    @(0.007287364918738604,0.9967244863510132,0.0002091415226459503,1),
    @(0.007142229005694389,0.9965912103652954,0.0003401002904865891,0.9997839331626892),
    @(0.007677559275180101,0.9963436126708984,0.0005145927425473928,0.9997839331626892),
    @(0.00825173407793045,0.9960960745811462,0.002199297538027167,0.9994598627090454),
    @(0.008022965863347054,0.9965150356292725,0.0002842725079972297,0.999891996383667),
    @(0.007770373020321131,0.9965340495109558,0.0003144124348182231,0.999891996383667),
    @(0.008068548515439034,0.9961341619491577,0.0003090054087806493,0.999891996383667),
    @(0.00816953182220459,0.9961341619491577,0.0003810987109318376,0.999891996383667),
    @(0.007525958120822906,0.996229350566864,0.0006407233886420727,0.9996759295463562),
    @(0.009212152101099491,0.9959818124771118,0.0005608850624412298,0.9997839331626892),
    @(0.008061543107032776,0.996400773525238,0.0005313611472956836,0.9997839331626892),
    @(0.007982369512319565,0.9962484240531921,0.0007055760943330824,0.9997839331626892),
    @(0.006867211312055588,0.9966863989830017,0.0004861875204369426,0.9997839331626892),
    @(0.008889189921319485,0.995753288269043,0.0002071004564641044,0.999891996383667),
    @(0.007957781665027142,0.9959246516227722,0.0009704586700536311,0.9994598627090454),
    @(0.007823970168828964,0.9964578747749329,0.001435461803339422,0.9994598627090454),
    @(0.007396026980131865,0.9966483116149902,0.0003451579250395298,0.9997839331626892),
    @(0.007174260448664427,0.9963817000389099,0.0004143759433645755,0.9997839331626892),
    @(0.00805114209651947,0.9960960745811462,0.00006273632607189938,1),
    @(0.007946528494358063,0.9960770010948181,0.0005030386964790523,0.9996759295463562),
    @(0.00748318200930953,0.9960579872131348,0.00009917288116412237,1),
    @(0.006906661204993725,0.9966863989830017,0.00005495289951795712,1),
    @(0.007994715124368668,0.9958294630050659,0.00003773679054575041,1)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Count; $j++) {
        $ws.Cells.Item($i + 2, $j + 1).Value = $row[$j]
    }
}

Write-Output "Updated training history values for 50 epochs"
